$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number: force Text format so
# Excel keeps them as strings (matching the source inlineStr cells)
# instead of auto-converting to numeric values, then clear the
# temporary format so no stray style survives on the cell.
$numericLookingRefs = @("D5", "D6", "D7", "D10", "D11", "D13", "D15", "D18", "D20", "D21", "D23", "D24", "D25", "D27", "D28", "D32", "D34", "D35", "D37", "D38", "D39", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50")
foreach ($r in $numericLookingRefs) {
    $ws.Range($r).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "64.446.10"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "2.628.34"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "596.63"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "152.82"
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  +5.03%  "
$ws.Range("D10").Value = "5.84"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("D11").Value = "0.394"
$ws.Range("E11").Value = "  +3.35%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "28.15"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").Value = "3.098.36"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").Value = "  +13.62%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "64.297.02"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "2.587.78"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "12.32"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").Value = "350.18"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "7.09"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").Value = "67.71"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").Value = "1.71"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "9.24"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "8.39"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").Value = "551.78"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D31").Value = "0.0₃0914"
$ws.Range("E31").Value = "  +8.45%  "
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("E33").Value = "  +4.99%  "
$ws.Range("D34").Value = "5.53"
$ws.Range("E34").Value = "  +5.60%  "
$ws.Range("D35").Value = "6.24"
$ws.Range("E35").Value = "  +2.38%  "
$ws.Range("E36").Value = "  +3.59%  "
$ws.Range("D37").Value = "166.10"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").Value = "20.09"
$ws.Range("E38").Value = "  +3.99%  "
$ws.Range("D39").Value = "2.00"
$ws.Range("E39").Value = "  +3.28%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "169.23"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("D43").Value = "4.11"
$ws.Range("E43").Value = "  +4.97%  "
$ws.Range("D44").Value = "23.22"
$ws.Range("E44").Value = "  +8.58%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.24"
$ws.Range("E45").Value = "  +13.15%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0592"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "0.642"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").Value = "0.0977"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("D50").Value = "19.39"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("E51").Value = "  +15.73%  "

# Strip the temporary Text format back off so the cells end up with no
# explicit style, same as the original workbook.
foreach ($r in $numericLookingRefs) {
    $ws.Range($r).ClearFormats()
}
